# edit.ps1 - apply "maj fiche de test" changes
# 1) "Modification des produis" (Request Title) -> "Création des produis"
# 2) "ModificationProduit.php" -> "CreationProduit.php"
# 3) "... en post modifie bien les données ..." -> "... en post insert bien les données ..."
#
# All remaining hunks in the source diff are purely structural
# (w:proofErr spell-check bracket removal / run re-splitting) and do not
# change any visible text, so they are not reproduced as separate edits.

$d = $word.ActiveDocument

# 1. Request Title cell: "Modification des " -> "Création des "
$d.Content.Find.Execute(
    "Modification des ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Création des ", 2
) | Out-Null

# 2. File name reference: "ModificationProduit.php" -> "CreationProduit.php"
$d.Content.Find.Execute(
    "ModificationProduit.php", $true, $false, $false, $false, $false,
    $true, 1, $false, "CreationProduit.php", 2
) | Out-Null

# 3. Verb change: "post modifie bien" -> "post insert bien"
$d.Content.Find.Execute(
    "post modifie bien", $true, $false, $false, $false, $false,
    $true, 1, $false, "post insert bien", 2
) | Out-Null

Write-Output "done"
